$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: the rule table is being trimmed from 5 columns (NAME / CONDITION x3 /
# ACTION) down to 3 (NAME / CONDITION / ACTION). The "ACTION" header slides
# left into C18 as the old D18/E18 header cells are dropped below.
$ws.Range("C18").Value2 = "ACTION"

# Row 19: likewise, the "Test" label slides left into C19.
$ws.Range("C19").Value2 = "Test"

# Columns D:E only ever held data for rows 18-25 of this rule table, and the
# table is being narrowed to just columns A:C - delete them outright so the
# sheet's dimension shrinks from A1:E25 to A1:C25.
$ws.Range("D18:E25").Delete()

# Row 25: the rule's label/name is updated.
$ws.Range("A25").Value2 = "Code changed 301236"
